$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36, pushing the existing rows 36-39 down to 37-40.
$ws.Rows.Item(36).Insert()

# Restore the standard row height used throughout the table.
$ws.Rows.Item(36).RowHeight = 13.05

# Populate the new row: a prospect ("Las Americas"-style row, no invoice date)
# for BEIRUT LEBANESE STREET FOOD, salesperson Norman, Ryan M (code 013).
$ws.Range("A36").Value = "BEIRUT LEBANESE STREET FOOD "
$ws.Range("B36").Value = "Norman, Ryan M"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("C36").Value = "013"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("E36").Value = "0008282"
$ws.Range("E36").NumberFormat = "@"
